$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 806.3125
$ws.Range("I2").Value = 254.6923
$ws.Range("K2").Value = 254.6923
$ws.Range("M2").Value = -141.6923
$ws.Range("H51").Value = 3498.5
$ws.Range("I51").Value = 3647
$ws.Range("J51").Value = 3350
$ws.Range("K51").Value = 3647
$ws.Range("L51").Value = 3350
$ws.Range("M51").Value = -3163
$ws.Range("N51").Value = -4318
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H64").Value = 5494.6665
$ws.Range("I64").Value = 5492
$ws.Range("K64").Value = 5492
$ws.Range("M64").Value = -5244
$ws.Range("H67").Value = 5494.6665
$ws.Range("I67").Value = 5492
$ws.Range("K67").Value = 5492
$ws.Range("M67").Value = -4634
$ws.Range("H69").Value = 48287.57
$ws.Range("J69").Value = 19666.666
$ws.Range("L69").Value = 58999.99800000001
$ws.Range("N69").Value = -60747.99800000001
$ws.Range("H72").Value = 48287.57
$ws.Range("J72").Value = 19666.666
$ws.Range("L72").Value = 176999.994
$ws.Range("N72").Value = -185735.994
$ws.Range("H132").Value = 829.3
$ws.Range("I132").Value = 829.3
$ws.Range("K132").Value = 2487.9
$ws.Range("M132").Value = 42.10000000000036
$ws.Range("H138").Value = 2538.9429
$ws.Range("I138").Value = 3181.0715
$ws.Range("J138").Value = 2110.8572
$ws.Range("K138").Value = 9543.2145
$ws.Range("L138").Value = 6332.571599999999
$ws.Range("M138").Value = -4403.2145
$ws.Range("N138").Value = -16612.5716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 39044.332
$ws.Range("J76").Value = 39044.332
$ws.Range("L76").Value = 39044.332
$ws.Range("N76").Value = -39720.332
$ws.Range("H79").Value = 39044.332
$ws.Range("J79").Value = 39044.332
$ws.Range("L79").Value = 39044.332
$ws.Range("N79").Value = -41384.332
$ws.Range("H96").Value = 44990
$ws.Range("J96").Value = 44990
$ws.Range("L96").Value = 44990
$ws.Range("N96").Value = -50482
$ws.Range("H97").Value = 949.5714
$ws.Range("I97").Value = 1099.4
$ws.Range("J97").Value = 575
$ws.Range("K97").Value = 1099.4
$ws.Range("L97").Value = 575
$ws.Range("M97").Value = -603.4000000000001
$ws.Range("N97").Value = -1567
$ws.Range("H124").Value = 85124.75
$ws.Range("J124").Value = 85124.75
$ws.Range("L124").Value = 85124.75
$ws.Range("N124").Value = -94944.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1274.125
$ws.Range("I94").Value = 1579
$ws.Range("J94").Value = 359.5
$ws.Range("K94").Value = 1579
$ws.Range("L94").Value = 359.5
$ws.Range("M94").Value = -1128
$ws.Range("N94").Value = -1261.5
$ws.Range("H99").Value = 1492.6
$ws.Range("I99").Value = 1492.6
$ws.Range("K99").Value = 1492.6
$ws.Range("M99").Value = 5.400000000000091
$ws.Range("H134").Value = 967.1429000000001
$ws.Range("I134").Value = 967.1429000000001
$ws.Range("K134").Value = 2901.4287
$ws.Range("M134").Value = -366.4287000000004

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1467.4667
$ws.Range("I31").Value = 1233.5
$ws.Range("J31").Value = 1734.8572
$ws.Range("K31").Value = 1233.5
$ws.Range("L31").Value = 1734.8572
$ws.Range("M31").Value = -938.5
$ws.Range("N31").Value = -2324.8572
$ws.Range("H34").Value = 1467.4667
$ws.Range("I34").Value = 1233.5
$ws.Range("J34").Value = 1734.8572
$ws.Range("K34").Value = 1233.5
$ws.Range("L34").Value = 1734.8572
$ws.Range("M34").Value = -1031.5
$ws.Range("N34").Value = -2138.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 384.5
$ws.Range("I11").Value = 323.875
$ws.Range("K11").Value = 971.625
$ws.Range("M11").Value = -831.625
$ws.Range("H12").Value = 7359.385
$ws.Range("J12").Value = 7968.5
$ws.Range("L12").Value = 23905.5
$ws.Range("N12").Value = -24251.5
$ws.Range("H23").Value = 150000140
$ws.Range("I23").Value = 300000000
$ws.Range("J23").Value = 299
$ws.Range("K23").Value = 900000000
$ws.Range("L23").Value = 897
$ws.Range("M23").Value = -899999765
$ws.Range("N23").Value = -1367
$ws.Range("H37").Value = 63987
$ws.Range("J37").Value = 63987
$ws.Range("L37").Value = 191961
$ws.Range("N37").Value = -192185
$ws.Range("H131").Value = 403082.7
$ws.Range("J131").Value = 403082.7
$ws.Range("L131").Value = 1209248.1
$ws.Range("N131").Value = -1219328.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1000000
$ws.Range("I7").Value = 1000000
$ws.Range("K7").Value = 1000000
$ws.Range("M7").Value = -999888
$ws.Range("H8").Value = 1000000
$ws.Range("I8").Value = 1000000
$ws.Range("K8").Value = 1000000
$ws.Range("M8").Value = -999861
$ws.Range("H80").Value = 2834
$ws.Range("I80").Value = 2332.6667
$ws.Range("J80").Value = 3335.3333
$ws.Range("K80").Value = 2332.6667
$ws.Range("L80").Value = 3335.3333
$ws.Range("M80").Value = -1334.6667
$ws.Range("N80").Value = -5331.3333
$ws.Range("H83").Value = 2834
$ws.Range("I83").Value = 2332.6667
$ws.Range("J83").Value = 3335.3333
$ws.Range("K83").Value = 11663.3335
$ws.Range("L83").Value = 16676.6665
$ws.Range("M83").Value = -6671.333500000001
$ws.Range("N83").Value = -26660.6665
$ws.Range("H107").Value = 2393.8333
$ws.Range("I107").Value = 380.83334
$ws.Range("J107").Value = 4406.8335
$ws.Range("K107").Value = 380.83334
$ws.Range("L107").Value = 4406.8335
$ws.Range("M107").Value = 1539.16666
$ws.Range("N107").Value = -8246.833500000001
$ws.Range("H126").Value = 3999.6
$ws.Range("I126").Value = 3999.5
$ws.Range("J126").Value = 3999.6667
$ws.Range("K126").Value = 11998.5
$ws.Range("L126").Value = 11999.0001
$ws.Range("M126").Value = -9528.5
$ws.Range("N126").Value = -16939.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4393.5
$ws.Range("I46").Value = 4393.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 4393.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -4205.5
$ws.Range("N46").ClearContents()
$ws.Range("H55").Value = 996.125
$ws.Range("I55").Value = 665.3
$ws.Range("K55").Value = 665.3
$ws.Range("M55").Value = -492.3
$ws.Range("H63").Value = 89077
$ws.Range("I63").Value = 89077
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 89077
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -88328
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 89077
$ws.Range("I66").Value = 89077
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 267231
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -263487
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 17141.5
$ws.Range("I62").Value = 24433.334
$ws.Range("J62").Value = 13495.583
$ws.Range("K62").Value = 24433.334
$ws.Range("L62").Value = 13495.583
$ws.Range("M62").Value = -23809.334
$ws.Range("N62").Value = -14743.583
$ws.Range("H65").Value = 17141.5
$ws.Range("I65").Value = 24433.334
$ws.Range("J65").Value = 13495.583
$ws.Range("K65").Value = 122166.67
$ws.Range("L65").Value = 67477.91500000001
$ws.Range("M65").Value = -119046.67
$ws.Range("N65").Value = -73717.91500000001
$ws.Range("H81").Value = 40002
$ws.Range("J81").Value = 40002
$ws.Range("L81").Value = 80004
$ws.Range("N81").Value = -82126
$ws.Range("H84").Value = 40002
$ws.Range("J84").Value = 40002
$ws.Range("L84").Value = 400020
$ws.Range("N84").Value = -410628
$ws.Range("H122").Value = 866.625
$ws.Range("I122").Value = 881.1429000000001
$ws.Range("J122").Value = 765
$ws.Range("K122").Value = 2643.4287
$ws.Range("L122").Value = 2295
$ws.Range("M122").Value = -193.4287000000004
$ws.Range("N122").Value = -7195
$ws.Range("H132").Value = 4611.6
$ws.Range("I132").Value = 5051.6665
$ws.Range("J132").Value = 3951.5
$ws.Range("K132").Value = 15154.9995
$ws.Range("L132").Value = 11854.5
$ws.Range("M132").Value = -12624.9995
$ws.Range("N132").Value = -16914.5
$ws.Range("H136").Value = 3493.1
$ws.Range("I136").Value = 3287.7856
$ws.Range("J136").Value = 3972.1667
$ws.Range("K136").Value = 9863.356800000001
$ws.Range("L136").Value = 11916.5001
$ws.Range("M136").Value = -7313.356800000001
$ws.Range("N136").Value = -17016.5001
